$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add header labels for new columns (X1:AH1), matching existing bold/border style ---
$headers = @(
    "wage_contr_to_ed",
    "cu_contr_to_ed",
    "ngdppot_contr_to_ed",
    "wage_contr_shortage",
    "cu_contr_shortage",
    "ngdppot_contr_shortage",
    "wage_contr_gcpi_via_ed",
    "cu_contr_gcpi_via_ed",
    "ngdppot_contr_gcpi_via_ed",
    "gcu_contr_gcpi",
    "gcu_contr_gw"
)

$startCol = 24  # column X
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# Copy the header style (bold font + border + centered alignment) from an existing header cell (W1)
# onto the new header range, so the new headers visually match the rest of row 1.
$ws.Range("W1").Copy() | Out-Null
$ws.Range("X1:AH1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Fill in the data values for rows 2-20, columns X-AH ---
$data = @{
    2 = @(-0.02922048654658926,-0.03428013632427085,0.03530281914597388,0,0,-0,0,0,-0,0,0)
    3 = @(-0.02257103463570331,-0.02148138566120839,0.02701215687793024,0,0,-0,0,0,-0,0,0)
    4 = @(-0.01596550603039404,-0.01174630630727247,0.01658896895558648,0,0,-0,0,0,-0,0,0)
    5 = @(-0.00722546496104215,-0.008013531667309337,0.008171313520501045,0,0,-0,0,0,-0,0,0)
    6 = @(0,-0,-0,0,0,0,0,0,0,0,0)
    7 = @(0.01002873734811605,0.01797784164150862,-0.009570804711598768,-0.4638241214031894,-0.8314662469114171,0.4426449644044898,-0.04981026362410087,-0.08929150305487896,0.04753582522221271,-0.005579484131193801,-0.01513906515331653)
    8 = @(0.01358619396364347,0.151029665444864,-0.01080228479886713,-0.4494422820657746,-4.996183454970099,0.3573483158358565,-0.04623461776667726,-0.5139628413931331,0.03676081101292281,0.02179755112064807,0.06700597961610644)
    9 = @(0.01925193149932092,0.0570863383509046,-0.02433943722991216,5.137671798238795,15.23436080273326,-6.495350362405663,0.5958176398652898,1.76673428253093,-0.7532681095653786,-0.002874316744263083,-0.03803185522194674)
    10 = @(0.02769080014518543,0.03380811693513852,-0.03596874897671043,-2.35185402389182,-2.871414167782195,3.05492244975591,-0.7673986580634803,-0.9369286344796877,0.996806504406806,0.6672117410865654,1.816429832215487)
    11 = @(0.03952525779218874,0.02560003899315788,-0.0534259920417739,-37.96857307754263,-24.59179283308175,51.32183308565772,-2.617753571642041,-1.695487828589405,3.538397705618817,-0.1480237479286686,-1.359058385987417)
    12 = @(0.04779546344085972,0.003527788320844216,-0.07372653861091294,-14.83048461656292,-1.09463967197447,22.87665435140111,-2.234846855327226,-0.1649542878657146,3.447346486610327,-0.1402262041891262,-0.2184911543699735)
    13 = @(0.06345802400381562,-0.005226721957527047,-0.0937641366886055,-3.629623999308326,0.298954084255937,5.363050081403819,1.252551093457823,-0.1031664065478445,-1.850741080897506,0.2238005034421153,0.3535042657350278)
    14 = @(0.07621415534825271,-0.01592334306286691,-0.1162555345731047,4.337734590342479,-0.9062783111887008,-6.616692808996318,-0.1038825724622378,0.02170407625736768,0.1584603796924925,-0.2038707602638858,-0.4700099198750243)
    15 = @(0.08880961520123165,-0.01922334998873865,-0.1417266661597019,0.5412848242629492,-0.1171642011602282,-0.8638084222277318,-0.3552937029762585,0.0769053574394144,0.5669948227594857,-0.1236228666022026,0.01198166855976179)
    16 = @(0.1031983526533304,-0.02303083280429785,-0.1696111612386186,-1.999661797231634,0.4462656169705559,3.286534627663119,0.325109545015457,-0.07255487497445724,-0.534332245060804,-0.01340588119028752,-0.06900073680243679)
    17 = @(0.1148188836763495,-0.01999482553647258,-0.1863278053582871,1.233603995280957,-0.2148226483046615,-2.001889565220194,-0.0384957408619013,0.006703737205812268,0.06247079470532135,-0.1372600118236607,-0.08632726983346561)
    18 = @(0.1263059276634806,-0.01104843406225964,-0.2011288848330306,1.857797933557207,-0.1625082713821334,-2.958347510158123,0.009788681531332758,-0.0008562512025843619,-0.01558744420632512,-0.06225814411308583,-0.02345471088259909)
    19 = @(0.1382896724646008,-0.007793395941804171,-0.2157638358458023,-1.650681587438936,0.09302513308107292,2.575444606371685,-0.01760139856652062,0.0009919371831147647,0.02746224792697688,-0.0235713039910781,0.004103230162844973)
    20 = @(0.1482710404226308,-0.00466377738107937,-0.2259159112090003,0.06588312639778197,-0.002072314551870309,-0.1003840432428832,0.0201731748905018,-0.0006345352166593454,-0.03073723077936015,-0.01525976423996189,0.08987924439261885)
}

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($row, $startCol + $j).Value = $rowValues[$j]
    }
}

Write-Output "Added columns X:AH with decomposition figures"